# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# 1) The "Periodo Mora" (col E) / "Valor Mora" (col F) block for rows 16-22
#    is reordered: it was listed newest-to-oldest (2311 down to 2305) and is
#    now listed oldest-to-newest (2305 up to 2311) - i.e. the 7-row block is
#    reversed in place.
# 2) The company logo picture is nudged to the left by 13.5pt (its size and
#    vertical position are unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Reverse the Periodo Mora / Valor Mora rows (E16:F22) ---------------

$periodos = @("2311", "2310", "2309", "2308", "2307", "2306", "2305")
$valores  = @(43307, 46400, 46400, 46400, 46400, 46400, 44028)

$firstRow = 16
$lastRow  = 22
$count    = $lastRow - $firstRow + 1

for ($i = 0; $i -lt $count; $i++) {
    $row = $firstRow + $i
    # reversed order: last item goes to the first row, and so on
    $srcIndex = $count - 1 - $i
    $ws.Cells.Item($row, 5).Value = $periodos[$srcIndex]
    $ws.Cells.Item($row, 6).Value = $valores[$srcIndex]
}

# --- 2) Reposition the logo picture -----------------------------------------

$shp = $ws.Shapes.Item(1)
$shp.Left = $shp.Left - 13.5
# The engine's Width/Height read-back drifts from the true EMU size once the
# shape has been touched, so pin them back to the picture's real size
# (975600 x 612000 EMU) instead of trusting a re-read of $shp.Width/Height.
$shp.Width = 76.81889763779527
$shp.Height = 48.188976377952756
